# chapter 14 finished, docs filled
# Fill in the log entries for 2020-05-23 (row 20), 2020-05-24 (row 21),
# and 2020-05-25 (row 22) on the time-tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 - 2020-05-23: 3 hours, chapter 13 finished / chapter 14 started reading,
# plus a short extra note about HTML intro.
# Set D20 before C20 so the shared-string table picks up the "short intro"
# text ahead of the "chapter 13" text, matching the source order.
$ws.Range("D20").Value = "short intro for HTML, good for revision"
$ws.Range("C20").Value = "chapter 13 completed, chapter 14 red"
$ws.Range("B20").Value = 3

# Row 21 - 2020-05-24: 0 hours, free day.
$ws.Range("B21").Value = 0
$ws.Range("D21").Value = "freeday"

# Row 22 - 2020-05-25: 3 hours, chapter 14 completed, docs filled, github updated.
$ws.Range("B22").Value = 3
$ws.Range("C22").Value = "chapter 14 completed, docs filled, github updated"

# Leave the selection on C22, matching the last-edited cell.
$ws.Range("C22").Select()
